# Remove the 10 blank spacer rows that separated the "z30" key/value block
# (rows 199-208, all empty) from the rest of the table. Everything below
# shifts up by 10 rows as a result (row 211 -> 201, ... row 352 -> 342),
# which matches the "modified key-val pairs for z30 in 32106" re-layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("199:208").Delete()

# Leave the view where the author left it after the edit.
$ws.Range("B47").Select()
